$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.811.39"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.740.78"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "231.03"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.5165"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("D8").Value = "0.2806"
$ws.Range("E8").Value = "  +5.06%  "
$ws.Range("D9").Value = "39.36"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").Value = "0.06104"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "1.770.09"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "0.07044"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "15.28"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "0.6411"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "4.515"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "77.05"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "25.829.40"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "0.000006579"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "1.975.25"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "4.130"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").Value = "8.634"
$ws.Range("E24").Value = "  +4.60%  "
$ws.Range("D25").Value = "5.148"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "139.83"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").Value = "1.505"
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("D28").Value = "15.08"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "1.817"
$ws.Range("E29").Value = "  +3.52%  "
$ws.Range("D30").Value = "102.60"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "0.08268"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").Value = "3.665"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").Value = "3.425"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "0.04490"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").Value = "2.613"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").Value = "0.9832"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").Value = "0.6153"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").Value = "2.652"
$ws.Range("E38").Value = "  +3.94%  "
$ws.Range("D39").Value = "0.01585"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "1.934"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "100.39"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "0.3836"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").Value = "0.7246"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").Value = "4.963"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").Value = "0.05411"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").Value = "6.275"
$ws.Range("E47").Value = "  +5.80%  "
$ws.Range("D48").Value = "0.1123"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").Value = "53.25"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "7.680"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").Value = "29.82"
$ws.Range("E51").Value = "  -0.71%  "
